$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates for the cryptocurrency price/volume refresh
$ws.Range('D2').Value = '26.893.76'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.667.83'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '215.33'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').Value = '0.522'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('D8').Value = '0.0623'
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('D10').Value = '20.32'
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('E11').Value = '  +3.17%  '
$ws.Range('D12').Value = '1.902.67'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').Value = '1.681.25'
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = '0.527'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '66.01'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '26.891.10'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '234.88'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '9.14'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '146.64'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '15.87'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').Value = '0.0495'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').Value = '1.445.92'
$ws.Range('E33').Value = '  -4.86%  '
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '0.585'
$ws.Range('E37').Value = '  +0.87%  '
$ws.Range('D38').Value = '0.904'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').Value = '5.74'
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '0.995'
$ws.Range('E42').Value = '  +8.85%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.29'
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').Value = '65.99'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '1.810.23'
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('D46').Value = '0.783'
$ws.Range('E46').Value = '  +1.12%  '
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('E50').Value = '  +4.09%  '
$ws.Range('E51').Value = '  +0.00%  '
